$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "202.173.124.126"
$ws.Range("B13").Value = 28.475392
$ws.Range("C13").Value = 77.0670592
$ws.Range("D13").Value = 616570.7228211587
$ws.Range("E13").Value = "Mozilla/5.0 (Windows NT 10.0; Win64; x64) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Safari/537.36"
$ws.Range("F13").Value = "Win32"
$ws.Range("G13").Value = "2025-06-25T16:52:19.204Z"
